# Sweden Superettan workbook update
# The commit swaps the data of several pairs of adjacent match rows (and one
# three-row rotation), while keeping the row-index column (A), Div (C),
# Div Original Name (D) and Date (E) columns untouched (they are identical
# between the affected rows anyway).  Columns B and F through AC hold the
# actual match data (id, teams, score, odds, ...) and are exchanged between
# the rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All data columns that participate in the swap (everything except A, C, D, E)
$cols = @('B','F','G','H','I','J','K','L','M','N','O','P','Q','R','S','T','U','V','W','X','Y','Z','AA','AB','AC')

function Swap-RowData {
    param($ws, [int]$row1, [int]$row2, $cols)

    foreach ($col in $cols) {
        $addr1 = "$col$row1"
        $addr2 = "$col$row2"
        $rng1 = $ws.Range($addr1)
        $rng2 = $ws.Range($addr2)
        $v1 = $rng1.Value2
        $v2 = $rng2.Value2
        if ($v1 -ne $v2) {
            $rng1.Value2 = $v2
            $rng2.Value2 = $v1
        }
    }
}

function Rotate-RowData {
    param($ws, [int[]]$rows, $cols)

    # new content of rows[i] becomes old content of rows[i+1] (wrap-around)
    $n = $rows.Length
    foreach ($col in $cols) {
        $old = @()
        foreach ($r in $rows) {
            $old += , ($ws.Range("$col$r").Value2)
        }
        for ($i = 0; $i -lt $n; $i++) {
            $newVal = $old[($i + 1) % $n]
            $rng = $ws.Range("$col$($rows[$i])")
            if ($rng.Value2 -ne $newVal) {
                $rng.Value2 = $newVal
            }
        }
    }
}

# Pairs of rows whose match data was swapped
$pairs = @(
    @(432, 433),
    @(484, 485),
    @(487, 488),
    @(549, 550),
    @(560, 561),
    @(570, 571),
    @(584, 585),
    @(588, 589),
    @(611, 612),
    @(666, 667),
    @(679, 680),
    @(695, 696)
)

foreach ($pair in $pairs) {
    Swap-RowData $ws $pair[0] $pair[1] $cols
}

# Three rows whose match data was cyclically rotated:
# new(661) = old(662), new(662) = old(663), new(663) = old(661)
Rotate-RowData $ws @(661, 662, 663) $cols
